$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 255
$ws.Range("C3").Value = 154332
$ws.Range("C4").Value = 145477
$ws.Range("C8").Value = 63.39
